# Project Sample Project2 is saved. Type: SAVE.
# The only data-level edit in this save is cell C8 on the active
# (Rules) sheet changing from 222 to 2222.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 2222
